# Re-style every table in the deck that still uses the old local table
# style ("Table_0", {F62A0091-EC77-4539-AE0B-138D4C2F5B4D}) so it uses the
# built-in style {CD9E9A76-01A4-4FEB-84C7-2B1DF6C89523} instead.
#
# Table.Style is read-only as a plain property assignment in this host
# ("Table styles cannot be assigned through a property - call
# Table.ApplyStyle("{GUID}") instead"), so we drive the change through
# Table.ApplyStyle(), scanning every slide/shape instead of hard-coding
# slide numbers so the script keeps working if the deck is reshuffled.

$oldStyleId = "{F62A0091-EC77-4539-AE0B-138D4C2F5B4D}"
$newStyleId = "{CD9E9A76-01A4-4FEB-84C7-2B1DF6C89523}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
